$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-15) down by one row (to 3-16) using Value2
# so numeric/date cells keep their raw serial values (Value2 avoids the
# locale-formatted-string round trip that .Value would otherwise perform).
$src = $ws.Range("A2:R15")
$dst = $ws.Range("A3:R16")
$dst.Value2 = $src.Value2

# Row 16's "Fecha" cell (D16) is brand-new territory beyond the original
# A1:R15 used range, so it starts out with the default (unstyled) format.
# Give it the same date-formatted style as the rest of column D.
$ws.Range("D2").Copy()
$ws.Range("D16").PasteSpecial(-4122)

# Write the brand-new week's record into row 2 (same constant columns as
# every other row in this subset: Mercado/Region/Categoria/etc.).
$ws.Range("A2").Value2 = 7
$ws.Range("B2").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value2 = "Ñuble"
$ws.Range("D2").Value2 = 44922
$ws.Range("E2").Value2 = 16
$ws.Range("F2").Value2 = 100114007
$ws.Range("G2").Value2 = "Jengibre"
$ws.Range("H2").Value2 = "Sin especificar"
$ws.Range("I2").Value2 = "Primera"
$ws.Range("J2").Value2 = 30
$ws.Range("K2").Value2 = 17000
$ws.Range("L2").Value2 = 17000
$ws.Range("M2").Value2 = 17000
$ws.Range("N2").Value2 = "`$/caja 13 kilos"
$ws.Range("O2").Value2 = "Perú"
$ws.Range("P2").Value2 = 1308
$ws.Range("Q2").Value2 = 13
$ws.Range("R2").Value2 = "Hortaliza"
